# Fixed issue with sample numbers 001-006 causing false positive data
# columns to be included in PCA plots.
#
# This updates the "UnitMass" values (column C) for the "+ loading" table
# (rows 2-21) and the "- loading" table (rows 23-42) on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- "+ loading" table (rows 2-21) ---
$ws.Range("C3").Value  = 57
$ws.Range("C4").Value  = 29
$ws.Range("C5").Value  = 175
$ws.Range("C6").Value  = 41
$ws.Range("C7").Value  = 91
$ws.Range("C8").Value  = 32
$ws.Range("C9").Value  = 103
$ws.Range("C10").Value = 115
$ws.Range("C11").Value = 55
$ws.Range("C12").Value = 56
$ws.Range("C13").Value = 40
$ws.Range("C14").Value = 231
$ws.Range("C15").Value = 53
$ws.Range("C16").Value = 112
$ws.Range("C17").Value = 67
$ws.Range("C18").Value = 119
$ws.Range("C19").Value = 51
$ws.Range("C20").Value = 128
$ws.Range("C21").Value = 216

# --- "- loading" table (rows 23-42) ---
$ws.Range("C24").Value = 23
$ws.Range("C25").Value = 58
$ws.Range("C26").Value = 97
$ws.Range("C27").Value = 102
$ws.Range("C29").Value = 24
$ws.Range("C30").Value = 85
$ws.Range("C31").Value = 30
$ws.Range("C32").Value = 125
$ws.Range("C33").Value = 74
$ws.Range("C34").Value = 70
$ws.Range("C35").Value = 28
$ws.Range("C37").Value = 60
$ws.Range("C38").Value = 98
$ws.Range("C39").Value = 86
$ws.Range("C40").Value = 82
$ws.Range("C41").Value = 72
$ws.Range("C42").Value = 138
